$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "USA" (sheet1): selection cursor moved to A5
# ---------------------------------------------------------------------
$usa = $wb.Worksheets.Item("USA")
$usa.Range("A5").Select()

# ---------------------------------------------------------------------
# Sheet "Canada" (sheet2): insert 3 new rows of data under "Central
# Canada" (row 6) for new satellite regions, then move selection.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Canada")

# Insert two blank rows after row 6 (old row 7 was already blank, so the
# net effect of the new rows 7-9 is a +2 shift for everything from the
# old row 8 onward). Use CopyOrigin = xlFormatFromLeftOrAbove (0) to
# avoid leaking the bordered/styled L column down into the blank rows.
$ws.Rows("7:8").Insert(-4121, 1)

# Clear any formatting the insert copied down into column L so the new
# rows stay completely empty there (matches target - no L7/L8/L9 cells).
$ws.Range("L7:L9").Clear()

# Row 6 gains new data in columns A/B (K6/L6 already existed).
$ws.Range("A6").Value = "Southern Ontario"
$ws.Range("B6").Value = '"http://www.ssd.noaa.gov/goes/east/gl/img/" + "YYYYDDD_hhmmoption.jpg"'
$ws.Range("B6").Characters(48, 7).Font.Color = 255
$ws.Range("B6").Characters(56, 10).Font.Color = 255

# New row 7: Southern Quebec
$ws.Range("A7").Value = "Southern Quebec"
$ws.Range("B7").Value = '"http://www.ssd.noaa.gov/goes/east/ne/img/ + "YYYYDDD_hhmmoption.jpg"'
$ws.Range("B7").Characters(47, 7).Font.Color = 255
$ws.Range("B7").Characters(55, 10).Font.Color = 255

# New row 8: Southern Western Canada
$ws.Range("A8").Value = "Southern Western Canada"
$ws.Range("B8").Value = '"http://www.ssd.noaa.gov/goes/west/nw/img/ + "YYYYDDD_hhmmoption.jpg"'
$ws.Range("B8").Characters(47, 7).Font.Color = 255
$ws.Range("B8").Characters(55, 10).Font.Color = 255

# New row 9: Southern Prairies
$ws.Range("A9").Value = "Southern Prairies"
$ws.Range("B9").Value = '"http://www.ssd.noaa.gov/goes/east/np/img/ + "YYYYDDD_hhmmoption.jpg"'
$ws.Range("B9").Characters(47, 7).Font.Color = 255
$ws.Range("B9").Characters(55, 10).Font.Color = 255

# Selection cursor moved to B9
$ws.Range("B9").Select()

# Page setup was touched (Print Setup dialog) leaving a portrait
# orientation page-setup record on the Canada sheet.
$ws.PageSetup.Orientation = 1
